# A new week of price data was inserted at row 12 ("Fruta / hortaliza, semanal"),
# pushing all the existing data rows (old 12..116) down by one row to (13..117).
# The sheet's used range grows from A1:R116 to A1:R117 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at position 12; every row from 12 downward
# (including the old row 12) shifts down by one.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new weekly record.
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(12, 3).Value = "La Araucanía"
$ws.Cells.Item(12, 4).Value = 45282
$ws.Cells.Item(12, 5).Value = 9
$ws.Cells.Item(12, 6).Value = 100112042
$ws.Cells.Item(12, 7).Value = "Locoto"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 40
$ws.Cells.Item(12, 11).Value = 2700
$ws.Cells.Item(12, 12).Value = 2700
$ws.Cells.Item(12, 13).Value = 2700
$ws.Cells.Item(12, 14).Value = "$/kilo"
$ws.Cells.Item(12, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(12, 16).Value = 2700
$ws.Cells.Item(12, 17).Value = 1
$ws.Cells.Item(12, 18).Value = "Hortaliza"
